# Fixed excel sheets creation
# Adds two new worksheets ("Dashboard" and "Emails") after the existing
# "Sheet" tab, each pre-populated with a header row.

$wb = $excel.ActiveWorkbook

# --- Dashboard sheet -------------------------------------------------
$wsDashboard = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsDashboard.Name = "Dashboard"
$wsDashboard.Range("A1").Value = "TotalNumber"
$wsDashboard.Range("B1").Value = "DomainsNumber"

# --- Emails sheet ------------------------------------------------------
$wsEmails = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsEmails.Name = "Emails"
$wsEmails.Range("A1").Value = "Email"
$wsEmails.Range("B1").Value = "Domain"
$wsEmails.Range("C1").Value = "DateAdded"

# Keep the original "Sheet" tab active/selected, as it was before the edit.
$wb.Worksheets.Item(1).Activate()
